$wb = $excel.ActiveWorkbook

# --- Sheet "Valeurs réelles": fill in values that became known this week,
#     and append the new week's row (date 2025-01-27 / serial 45684) ---
$ws1 = $wb.Worksheets.Item("Valeurs réelles")

$ws1.Range("E23").Value = 19
$ws1.Range("D24").Value = 19
$ws1.Range("C25").Value = 19

$ws1.Range("A26").Value = 45684
$ws1.Range("A26").NumberFormat = $ws1.Range("A25").NumberFormat
$ws1.Range("B26").Value = 19

# keep C26/D26/E26 present as (still) empty cells, matching the pattern of
# the other not-yet-known observations on this sheet
$ws1.Range("B2").Copy($ws1.Range("C26"))
$ws1.Range("B2").Copy($ws1.Range("D26"))
$ws1.Range("B2").Copy($ws1.Range("E26"))

# --- Sheet "Prédictions": append the new week's predictions row ---
$ws2 = $wb.Worksheets.Item("Prédictions")

$ws2.Range("A26").Value = 45684
$ws2.Range("A26").NumberFormat = $ws2.Range("A25").NumberFormat
$ws2.Range("B26").Value = 15.34035682678223
$ws2.Range("C26").Value = 16.72483825683594
$ws2.Range("D26").Value = 16.46977424621582
